$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-37: date serial 45657 -> 45658
for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 3).Value = 45658
}

# Swap row 36 and 37 values for column A (Beteckning) and column G (Area)
$ws.Range("A36").Value = "A 60500-2024"
$ws.Range("G36").Value = 0.8

$ws.Range("A37").Value = "A 60501-2024"
$ws.Range("G37").Value = 0.6
